# Auto-generated Excel COM-interop script to update column F ("想去人数")
# values across sheets 展览, 演出, 全部类型 per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(6, 6).Value = 306
$ws.Cells.Item(8, 6).Value = 11
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(10, 6).Value = 499
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(13, 6).Value = 36
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(16, 6).Value = 6
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(22, 6).Value = 676
$ws.Cells.Item(23, 6).Value = 46
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(27, 6).Value = 1082
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(30, 6).Value = 175
$ws.Cells.Item(31, 6).Value = 5219
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(37, 6).Value = 31
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(41, 6).Value = 20
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(45, 6).Value = 3983
$ws.Cells.Item(46, 6).Value = 0

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(3, 6).Value = 0

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(7, 6).Value = 1092
$ws.Cells.Item(8, 6).Value = 11
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(12, 6).Value = 257
$ws.Cells.Item(13, 6).Value = 36
$ws.Cells.Item(14, 6).Value = 151
$ws.Cells.Item(16, 6).Value = 6
$ws.Cells.Item(17, 6).Value = 232
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(22, 6).Value = 676
$ws.Cells.Item(23, 6).Value = 46
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(31, 6).Value = 5220
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(35, 6).Value = 33
$ws.Cells.Item(36, 6).Value = 2796
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(39, 6).Value = 31
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(43, 6).Value = 20
$ws.Cells.Item(45, 6).Value = 254
$ws.Cells.Item(47, 6).Value = 3983
$ws.Cells.Item(48, 6).Value = 0
$ws.Cells.Item(49, 6).Value = 93
